$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Probabilities")

$ws.Cells.Item(2,1).Value = 1369
$ws.Cells.Item(2,2).Value = "2025-11-28T17:00:00"
$ws.Cells.Item(2,3).Value = "Трактор"
$ws.Cells.Item(2,4).Value = "ХК Сочи"
$ws.Cells.Item(2,5).Value = 897827
$ws.Cells.Item(2,6).Value = "https://text.khl.ru/text/897827.html"
$ws.Cells.Item(2,7).Value = 5.5
$ws.Cells.Item(2,8).Value = 1
$ws.Cells.Item(2,9).Value = 4.7
$ws.Cells.Item(2,10).Value = 6.875
$ws.Cells.Item(2,11).Value = 6.1875
$ws.Cells.Item(2,12).Value = 2.85
$ws.Cells.Item(2,13).Value = 6.5
$ws.Cells.Item(2,14).Value = 41.132791
$ws.Cells.Item(2,15).Value = 22.883368
$ws.Cells.Item(2,16).Value = 64.016159
$ws.Cells.Item(2,17).Value = 0.2
$ws.Cells.Item(2,18).Value = -0.2
$ws.Cells.Item(2,19).Value = 0.779471
$ws.Cells.Item(2,20).Value = 0.073832
$ws.Cells.Item(2,21).Value = 0.09566
$ws.Cells.Item(2,22).Value = 0.020671
$ws.Cells.Item(2,23).Value = 0.928292
$ws.Cells.Item(2,24).Value = 0.053712
$ws.Cells.Item(2,25).Value = 0.895251
$ws.Cells.Item(2,26).Value = 0.113432
$ws.Cells.Item(2,27).Value = 0.835531
$ws.Cells.Item(2,28).Value = 0.203386
$ws.Cells.Item(2,29).Value = 0.745577
$ws.Cells.Item(2,30).Value = 0.319523
$ws.Cells.Item(2,31).Value = 0.62944
$ws.Cells.Item(2,32).Value = 0.98523
$ws.Cells.Item(2,33).Value = 0.01477
$ws.Cells.Item(2,34).Value = 0.945893
$ws.Cells.Item(2,35).Value = 0.054107
$ws.Cells.Item(2,36).Value = 0.777299
$ws.Cells.Item(2,37).Value = 0.222701
$ws.Cells.Item(2,38).Value = 0.542379
$ws.Cells.Item(2,39).Value = 0.457621
$ws.Cells.Item(2,40).Value = 0.90032
$ws.Cells.Item(2,41).Value = 0.271571

$ws.Cells.Item(3,1).Value = 1369
$ws.Cells.Item(3,2).Value = "2025-11-28T19:00:00"
$ws.Cells.Item(3,3).Value = "Ак Барс"
$ws.Cells.Item(3,4).Value = "СКА"
$ws.Cells.Item(3,5).Value = 897823
$ws.Cells.Item(3,6).Value = "https://text.khl.ru/text/897823.html"
$ws.Cells.Item(3,7).Value = 2.240677
$ws.Cells.Item(3,8).Value = 3.923077
$ws.Cells.Item(3,9).Value = 4.6875
$ws.Cells.Item(3,10).Value = 1.756236
$ws.Cells.Item(3,11).Value = 1.998456
$ws.Cells.Item(3,12).Value = 4.305288
$ws.Cells.Item(3,13).Value = 6.163754
$ws.Cells.Item(3,14).Value = 29.335543
$ws.Cells.Item(3,15).Value = 34.435929
$ws.Cells.Item(3,16).Value = 63.771472
$ws.Cells.Item(3,17).Value = -0.101243
$ws.Cells.Item(3,18).Value = 0.2
$ws.Cells.Item(3,19).Value = 0.125662
$ws.Cells.Item(3,20).Value = 0.108911
$ws.Cells.Item(3,21).Value = 0.760551
$ws.Cells.Item(3,22).Value = 0.126087
$ws.Cells.Item(3,23).Value = 0.869037
$ws.Cells.Item(3,24).Value = 0.246453
$ws.Cells.Item(3,25).Value = 0.748671
$ws.Cells.Item(3,26).Value = 0.398203
$ws.Cells.Item(3,27).Value = 0.596921
$ws.Cells.Item(3,28).Value = 0.557636
$ws.Cells.Item(3,29).Value = 0.437488
$ws.Cells.Item(3,30).Value = 0.701211
$ws.Cells.Item(3,31).Value = 0.293913
$ws.Cells.Item(3,32).Value = 0.593576
$ws.Cells.Item(3,33).Value = 0.406424
$ws.Cells.Item(3,34).Value = 0.322906
$ws.Cells.Item(3,35).Value = 0.677094
$ws.Cells.Item(3,36).Value = 0.928395
$ws.Cells.Item(3,37).Value = 0.071605
$ws.Cells.Item(3,38).Value = 0.803308
$ws.Cells.Item(3,39).Value = 0.196692
$ws.Cells.Item(3,40).Value = 0.380078
$ws.Cells.Item(3,41).Value = 0.937003

$ws.Cells.Item(4,1).Value = 1369
$ws.Cells.Item(4,2).Value = "2025-11-28T19:00:00"
$ws.Cells.Item(4,3).Value = "Торпедо"
$ws.Cells.Item(4,4).Value = "Динамо Мн"
$ws.Cells.Item(4,5).Value = 897824
$ws.Cells.Item(4,6).Value = "https://text.khl.ru/text/897824.html"
$ws.Cells.Item(4,7).Value = 2.383937
$ws.Cells.Item(4,8).Value = 4.636364
$ws.Cells.Item(4,9).Value = 3.009599
$ws.Cells.Item(4,10).Value = 1.484589
$ws.Cells.Item(4,11).Value = 1.934263
$ws.Cells.Item(4,12).Value = 3.822981
$ws.Cells.Item(4,13).Value = 7.0203
$ws.Cells.Item(4,14).Value = 27.808605
$ws.Cells.Item(4,15).Value = 39.608055
$ws.Cells.Item(4,16).Value = 67.41666
$ws.Cells.Item(4,17).Value = -0.060951
$ws.Cells.Item(4,18).Value = 0.2
$ws.Cells.Item(4,19).Value = 0.154026
$ws.Cells.Item(4,20).Value = 0.127631
$ws.Cells.Item(4,21).Value = 0.716317
$ws.Cells.Item(4,22).Value = 0.174216
$ws.Cells.Item(4,23).Value = 0.823758
$ws.Cells.Item(4,24).Value = 0.318863
$ws.Cells.Item(4,25).Value = 0.679112
$ws.Cells.Item(4,26).Value = 0.485416
$ws.Cells.Item(4,27).Value = 0.512559
$ws.Cells.Item(4,28).Value = 0.64523
$ws.Cells.Item(4,29).Value = 0.352744
$ws.Cells.Item(4,30).Value = 0.776672
$ws.Cells.Item(4,31).Value = 0.221303
$ws.Cells.Item(4,32).Value = 0.575909
$ws.Cells.Item(4,33).Value = 0.424091
$ws.Cells.Item(4,34).Value = 0.305537
$ws.Cells.Item(4,35).Value = 0.694463
$ws.Cells.Item(4,36).Value = 0.894557
$ws.Cells.Item(4,37).Value = 0.105443
$ws.Cells.Item(4,38).Value = 0.734795
$ws.Cells.Item(4,39).Value = 0.265205
$ws.Cells.Item(4,40).Value = 0.443622
$ws.Cells.Item(4,41).Value = 0.925895

$ws.Cells.Item(5,1).Value = 1369
$ws.Cells.Item(5,2).Value = "2025-11-28T19:00:00"
$ws.Cells.Item(5,3).Value = "Северсталь"
$ws.Cells.Item(5,4).Value = "Локомотив"
$ws.Cells.Item(5,5).Value = 897825
$ws.Cells.Item(5,6).Value = "https://text.khl.ru/text/897825.html"
$ws.Cells.Item(5,7).Value = 1.5
$ws.Cells.Item(5,8).Value = 1.40625
$ws.Cells.Item(5,9).Value = 1.233333
$ws.Cells.Item(5,10).Value = 1.15625
$ws.Cells.Item(5,11).Value = 1.328125
$ws.Cells.Item(5,12).Value = 1.319792
$ws.Cells.Item(5,13).Value = 2.90625
$ws.Cells.Item(5,14).Value = 22.31792
$ws.Cells.Item(5,15).Value = 22.352212
$ws.Cells.Item(5,16).Value = 44.670133
$ws.Cells.Item(5,17).Value = -0.2
$ws.Cells.Item(5,18).Value = -0.2
$ws.Cells.Item(5,19).Value = 0.371378
$ws.Cells.Item(5,20).Value = 0.261117
$ws.Cells.Item(5,21).Value = 0.367505
$ws.Cells.Item(5,22).Value = 0.72554
$ws.Cells.Item(5,23).Value = 0.27446
$ws.Cells.Item(5,24).Value = 0.87056
$ws.Cells.Item(5,25).Value = 0.129439
$ws.Cells.Item(5,26).Value = 0.947361
$ws.Cells.Item(5,27).Value = 0.052639
$ws.Cells.Item(5,28).Value = 0.981255
$ws.Cells.Item(5,29).Value = 0.018745
$ws.Cells.Item(5,30).Value = 0.994076
$ws.Cells.Item(5,31).Value = 0.005924
$ws.Cells.Item(5,32).Value = 0.383108
$ws.Cells.Item(5,33).Value = 0.616892
$ws.Cells.Item(5,34).Value = 0.149413
$ws.Cells.Item(5,35).Value = 0.850587
$ws.Cells.Item(5,36).Value = 0.380173
$ws.Cells.Item(5,37).Value = 0.619827
$ws.Cells.Item(5,38).Value = 0.147469
$ws.Cells.Item(5,39).Value = 0.852531
$ws.Cells.Item(5,40).Value = 0.835571
$ws.Cells.Item(5,41).Value = 0.832979

$ws.Cells.Item(6,1).Value = 1369
$ws.Cells.Item(6,2).Value = "2025-11-28T19:00:00"
$ws.Cells.Item(6,3).Value = "Нефтехимик"
$ws.Cells.Item(6,4).Value = "Драконы"
$ws.Cells.Item(6,5).Value = 897826
$ws.Cells.Item(6,6).Value = "https://text.khl.ru/text/897826.html"
$ws.Cells.Item(6,7).Value = 1.591473
$ws.Cells.Item(6,8).Value = 1.848538
$ws.Cells.Item(6,9).Value = 4.736842
$ws.Cells.Item(6,10).Value = 5.205882
$ws.Cells.Item(6,11).Value = 3.398678
$ws.Cells.Item(6,12).Value = 3.29269
$ws.Cells.Item(6,13).Value = 3.440011
$ws.Cells.Item(6,14).Value = 27.688566
$ws.Cells.Item(6,15).Value = 27.29634
$ws.Cells.Item(6,16).Value = 54.984906
$ws.Cells.Item(6,17).Value = -0.142655
$ws.Cells.Item(6,18).Value = -0.107673
$ws.Cells.Item(6,19).Value = 0.436582
$ws.Cells.Item(6,20).Value = 0.157276
$ws.Cells.Item(6,21).Value = 0.404708
$ws.Cells.Item(6,22).Value = 0.099342
$ws.Cells.Item(6,23).Value = 0.899224
$ws.Cells.Item(6,24).Value = 0.203053
$ws.Cells.Item(6,25).Value = 0.795514
$ws.Cells.Item(6,26).Value = 0.341846
$ws.Cells.Item(6,27).Value = 0.65672
$ws.Cells.Item(6,28).Value = 0.496633
$ws.Cells.Item(6,29).Value = 0.501934
$ws.Cells.Item(6,30).Value = 0.644594
$ws.Cells.Item(6,31).Value = 0.353972
$ws.Cells.Item(6,32).Value = 0.853008
$ws.Cells.Item(6,33).Value = 0.146992
$ws.Cells.Item(6,34).Value = 0.660005
$ws.Cells.Item(6,35).Value = 0.339995
$ws.Cells.Item(6,36).Value = 0.84051
$ws.Cells.Item(6,37).Value = 0.15949
$ws.Cells.Item(6,38).Value = 0.639103
$ws.Cells.Item(6,39).Value = 0.360897
$ws.Cells.Item(6,40).Value = 0.736572
$ws.Cells.Item(6,41).Value = 0.709293
